$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "timestamp" column (O) for all data rows (2-73) from the
# old scrape time to the new one.
$ws.Range("O2:O73").Value = "2022-08-04 20:57:44"
